# Daily "cryptos" list refresh (GitHub Actions bot edit).
# Updates Price (D) / Volume(1h) (E) columns for each coin row, and swaps
# the Stacks/Maker rows (42/43) to reflect the new ranking order.
#
# Note: several Price values look like plain decimals (e.g. "54.90"), which
# Excel's cell-value parser would otherwise auto-coerce to a Number. Those
# are written with a leading apostrophe (Excel's native "force text" quote
# prefix) so they round-trip as text, matching the source feed's formatting
# (thousand-separator dots elsewhere, e.g. "70.866.21", make that intent
# clear). Values that already contain a second "." or other non-numeric
# character stay naturally text and are written plain.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "70.866.21"
$ws.Range("E2").Value = "  +1.77%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.589.51"
$ws.Range("E3").Value = "  +1.20%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'590.74"
$ws.Range("E5").Value = "  +2.81%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'187.68"
$ws.Range("E6").Value = "  +1.47%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.579.23"
$ws.Range("E7").Value = "  +1.11%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.44%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.05%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.203"
$ws.Range("E10").Value = "  +12.20%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "'0.654"
$ws.Range("E11").Value = "  +1.76%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'54.90"
$ws.Range("E12").Value = "  +1.05%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000313"
$ws.Range("E13").Value = "  +5.57%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +2.11%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.150.51"
$ws.Range("E15").Value = "  +0.68%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'19.50"
$ws.Range("E16").Value = "  +0.52%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "70.754.94"
$ws.Range("E17").Value = "  +1.76%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.592.12"
$ws.Range("E18").Value = "  +0.95%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'12.51"
$ws.Range("E19").Value = "  +0.95%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'566.72"
$ws.Range("E20").Value = "  +17.19%  "

# Row 21 - TRON
$ws.Range("E21").Value = "  +0.08%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'1.03"
$ws.Range("E22").Value = "  +0.64%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "'17.90"
$ws.Range("E23").Value = "  -7.03%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  +8.87%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +2.26%  "

# Row 26 - Litecoin
$ws.Range("D26").Value = "'96.22"
$ws.Range("E26").Value = "  +1.59%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "'11.64"
$ws.Range("E27").Value = "  +2.71%  "

# Row 28 - ImmutableX
$ws.Range("E28").Value = "  +3.09%  "

# Row 29 - Filecoin
$ws.Range("D29").Value = "'9.21"
$ws.Range("E29").Value = "  +0.31%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "'32.45"
$ws.Range("E30").Value = "  +3.70%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'7.36"
$ws.Range("E31").Value = "  -0.85%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "'12.61"
$ws.Range("E32").Value = "  +6.02%  "

# Row 33 - OKB
$ws.Range("D33").Value = "'65.28"
$ws.Range("E33").Value = "  -1.62%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +2.56%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "'563.39"
$ws.Range("E35").Value = "  -0.18%  "

# Row 36 - Fetch.AI
$ws.Range("D36").Value = "'3.29"
$ws.Range("E36").Value = "  +4.48%  "

# Row 37 - TheGraph
$ws.Range("E37").Value = "  +7.10%  "

# Row 38 - InjectiveProtocol
$ws.Range("D38").Value = "'38.32"
$ws.Range("E38").Value = "  +0.24%  "

# Row 39 - Dai
$ws.Range("E39").Value = "  +0.10%  "

# Row 40 - PEPE
$ws.Range("D40").Value = "0.0₃0779"
$ws.Range("E40").Value = "  -0.07%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +1.47%  "

# Rows 42/43 - ranking swap: Maker moves up to 42, Stacks drops to 43
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.363.09"
$ws.Range("E42").Value = "  +5.26%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'3.38"
$ws.Range("E43").Value = "  -2.31%  "

# Row 44 - dogwifhat
$ws.Range("D44").Value = "'3.07"
$ws.Range("E44").Value = "  -1.90%  "

# Row 45 - ApeXProtocol
$ws.Range("D45").Value = "'3.58"
$ws.Range("E45").Value = "  +4.41%  "

# Row 46 - ThetaToken
$ws.Range("E46").Value = "  +1.13%  "

# Row 47 - VeChain
$ws.Range("D47").Value = "'0.0448"
$ws.Range("E47").Value = "  +3.49%  "

# Row 48 - THORChain
$ws.Range("D48").Value = "'9.43"
$ws.Range("E48").Value = "  +0.09%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  +2.30%  "

# Row 50 - FirstDigitalUSD
$ws.Range("D50").Value = "'0.998"
$ws.Range("E50").Value = "  +0.17%  "

# Row 51 - OceanProtocol
$ws.Range("E51").Value = "  +20.15%  "
